# Update "想去人数" (want-to-go count) values in column F for rows 2,3,4,5,11,12,13,14
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 162
    3  = 7134
    4  = 5142
    5  = 75
    11 = 84
    12 = 192
    13 = 629
    14 = 190
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
